# Fix udi's display and score logic
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Update F6 balance from 1000 to 1002
$ws.Range("F6").Value = 1002

# Append 4 more rows (56-59), duplicating the data already present in row 55
# (moses / bro / 1234 / m@g.c / Male / balance 0) by copying the row so the
# underlying cell types (shared-string text) and styles are preserved exactly.
for ($r = 56; $r -le 59; $r++) {
    $ws.Range("A55:F55").Copy()
    $ws.Range("A" + $r + ":F" + $r).PasteSpecial(-4104)
}
$excel.CutCopyMode = $false
